$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to reflect repulled / recalculated data
$ws.Range("F2").Value = 4
$ws.Range("F5").Value = -7
$ws.Range("F9").Value = -2
$ws.Range("F13").Value = -6
$ws.Range("F15").Value = -3
$ws.Range("F17").Value = -3
$ws.Range("F19").Value = -6
$ws.Range("F20").Value = 1
$ws.Range("F23").Value = 2
$ws.Range("F25").Value = 4
$ws.Range("F26").Value = -4
$ws.Range("F27").Value = -5
